$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: Config Name -> "S1", End Time -> 44.72666666666667
$ws.Range("A3").Value = "S1"
$ws.Range("I3").Value = 44.72666666666667

# Update row 4: Config Name -> "Test", Plateau Time -> 44.73, End Time -> 3
$ws.Range("A4").Value = "Test"
$ws.Range("D4").Value = 44.73
$ws.Range("I4").Value = 3
